$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new observation was inserted as row 230 in the data table (which starts
# at row 2, header at row 1). This pushes all subsequent rows (old 230-297)
# down by one (new 231-298), and the sheet's used range grows from
# A1:R297 to A1:R298.

# Insert a new row above current row 230 (shifts old 230..297 down to 231..298,
# copying formatting/formulas from the row above as Excel normally does not,
# but we will explicitly set style/values for row 230 below).
$ws.Rows.Item(230).Insert()

# Populate the newly inserted row 230 with the new data point.
$ws.Cells.Item(230, 1).Value = 4
$ws.Cells.Item(230, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(230, 3).Value = "Los Lagos"
$ws.Cells.Item(230, 4).Value = 44841
$ws.Cells.Item(230, 4).NumberFormat = $ws.Cells.Item(231, 4).NumberFormat
$ws.Cells.Item(230, 5).Value = 10
$ws.Cells.Item(230, 6).Value = 100112044
$ws.Cells.Item(230, 7).Value = "Perejil"
$ws.Cells.Item(230, 8).Value = "Sin especificar"
$ws.Cells.Item(230, 9).Value = "Primera"
$ws.Cells.Item(230, 10).Value = 160
$ws.Cells.Item(230, 11).Value = 5000
$ws.Cells.Item(230, 12).Value = 5000
$ws.Cells.Item(230, 13).Value = 5000
$ws.Cells.Item(230, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(230, 15).Value = "Región Metropolitana"
$ws.Cells.Item(230, 16).Value = 1667
$ws.Cells.Item(230, 17).Value = 3
$ws.Cells.Item(230, 18).Value = "Hortaliza"
